$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -66.666666666666
$ws.Range("I16").Value = 2
$ws.Range("K16").Value = -66.666666666666
$ws.Range("L16").Value = -80
$ws.Range("M16").Value = -60
$ws.Range("N16").Value = -95.833333333333
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 2
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -83.333333333333
$ws.Range("I17").Value = 2
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = -83.333333333333
$ws.Range("L17").Value = -66.666666666666
$ws.Range("M17").Value = -71.428571428571
$ws.Range("N17").Value = -87.5
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 10
$ws.Range("K18").Value = -20
$ws.Range("L18").Value = -42.857142857142
$ws.Range("N18").Value = -90.804597701149
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 22.222222222222
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -29.090909090909
$ws.Range("I19").Value = 36
$ws.Range("J19").Value = 48
$ws.Range("K19").Value = -25
$ws.Range("L19").Value = 5.882352941176
$ws.Range("M19").Value = -23.404255319148
$ws.Range("N19").Value = -72.093023255813
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 3
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = -25
$ws.Range("N20").Value = -93.181818181818
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -5.882352941176
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 86
$ws.Range("H21").Value = -36.046511627907
$ws.Range("I21").Value = 52
$ws.Range("J21").Value = 79
$ws.Range("K21").Value = -34.177215189873
$ws.Range("L21").Value = -23.529411764705
$ws.Range("M21").Value = -22.388059701492
$ws.Range("N21").Value = -84
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = -50
$ws.Range("M22").Value = 0
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 53.968253968254
$ws.Range("I24").Value = 94
$ws.Range("J24").Value = 58
$ws.Range("K24").Value = 62.068965517241
$ws.Range("L24").Value = 42.424242424242
$ws.Range("M24").Value = 108.888888888889
$ws.Range("C25").Value = 25
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 177.777777777778
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 43.75
$ws.Range("I25").Value = 69
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = 60.465116279069
$ws.Range("L25").Value = 40.816326530612
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = -25
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -42.105263157894
$ws.Range("I26").Value = 7
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = -63.157894736842
$ws.Range("L26").Value = -65
$ws.Range("M26").Value = -63.157894736842
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = -50
$ws.Range("J28").Value = 2
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -66.666666666666

# --- Numeric cells becoming text placeholders ("0" / "***.*") ---
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("A16").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("A16").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("A20").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("A20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("A27").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("A31").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("A31").Copy()
$ws.Range("E31").PasteSpecial(-4122)

# --- Text placeholder cells becoming numeric ---
$ws.Range("C16").Value = 1
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C17").Value = 1
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("C22").Value = 2
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 2
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("I22").Value = 2
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("L22").Value = 0
$ws.Range("L22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C26").Value = 3
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"

$excel.CutCopyMode = $false
